# 2018_termination_data.xlsx - add monthly termination report
#
# This mirrors the authoring change described in the commit:
#   "add monthly data to reports ..."
# A new worksheet "2018_monthly" is inserted just before the
# "2018_procedure_location_wrksht" worksheet and populated with the
# Indiana 2018 resident / non-resident termination counts by month.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2018_monthly" worksheet right before the existing
#    "2018_procedure_location_wrksht" sheet (this keeps the relationship
#    ids / sheetIds of the existing sheets the same way Excel would when
#    you right click a tab -> Insert -> Worksheet).
# ---------------------------------------------------------------------
$wrksht = $wb.Worksheets.Item("2018_procedure_location_wrksht")
$monthly = $wb.Worksheets.Add($wrksht)
$monthly.Name = "2018_monthly"

# ---------------------------------------------------------------------
# 2. Fill in the monthly resident / non-resident termination data.
#    (Note: literal CR characters are used inside the header cells,
#    matching the original author's multi-line header text.)
# ---------------------------------------------------------------------
$CR = [char]13

$headerResident = "Resident" + $CR + "Terminations" + $CR + "(n =7263)"
$headerNonResident = "Non-Resident" + $CR + "Terminations" + $CR + "(n =774)"

$monthlyData = @(
    @("Month", $headerResident, $headerNonResident),
    @("January",   "576 (7.17%)", "59 (.73%)"),
    @("February",  "626 (7.79%)", "40 (.50%)"),
    @("March",     "746 (9.28%)", "78 (.97%)"),
    @("April",     "581 (7.23%)", "65 (.81%)"),
    @("May",       "666 (8.29%)", "72 (.89%)"),
    @("June",      "641 (7.98%)", "74 (.92%)"),
    @("July",      "494 (6.15%)", "58 (.72%)"),
    @("August",    "668 (8.31%)", "76 (.94%)"),
    @("September", "579 (7.20%)", "57 (.70%)"),
    @("October",   "523 (6.51%)", "63 (.78%)"),
    @("November",  "588 (7.32%)", "66 (.82%)"),
    @("December",  "575 (7.15%)", "66 (.82%)")
)

for ($r = 0; $r -lt $monthlyData.Length; $r++) {
    $rowValues = $monthlyData[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $monthly.Cells.Item($r + 1, $c + 1).Value = $rowValues[$c]
    }
}

# Match the author's final selection/view on the new sheet and make it
# the active tab (this is the sheet that was being edited when saved).
$monthly.Range("E12").Select()
$monthly.Activate()

# ---------------------------------------------------------------------
# 3. Minor cleanup that happened on the existing report sheets when the
#    workbook was re-saved: a handful of percent cells had a redundant
#    "General" number format explicitly applied - clear that back to
#    the plain default format.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2018_ethnicity").Range("C2:C4").ClearFormats()
$wb.Worksheets.Item("2018_education").Range("C2:C10").ClearFormats()
$wb.Worksheets.Item("2018_provider_location").Range("D2:D11").ClearFormats()
$wb.Worksheets.Item("2018_gestation").Range("C2:C5").ClearFormats()
